$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.302699565887451
$ws.Range("B1").Value = 4.221007347106934
$ws.Range("C1").Value = 2.712152242660522
$ws.Range("D1").Value = 1.889561295509338
$ws.Range("E1").Value = 1.199711322784424
